$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed by Excel as a
# number (single decimal point, no thousands separators) need the
# cell pre-formatted as Text so the literal string is preserved,
# matching the original inlineStr/Text cell content exactly.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply the updated values (coin names, links, prices, 1h volume %).
$ws.Range("D2").Value = '60.791.56'
$ws.Range("E2").Value = '  +3.02%  '
$ws.Range("D3").Value = '2.692.81'
$ws.Range("E3").Value = '  +2.06%  '
$ws.Range("D5").Value = '525.81'
$ws.Range("E5").Value = '  +1.75%  '
$ws.Range("D6").Value = '145.42'
$ws.Range("E6").Value = '  +1.16%  '
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  +1.97%  '
$ws.Range("D9").Value = '2.716.75'
$ws.Range("E9").Value = '  +1.91%  '
$ws.Range("D10").Value = '6.53'
$ws.Range("E10").Value = '  +4.68%  '
$ws.Range("E11").Value = '  +1.20%  '
$ws.Range("D12").Value = '0.340'
$ws.Range("E12").Value = '  +0.93%  '
$ws.Range("E13").Value = '  +2.62%  '
$ws.Range("D14").Value = '3.167.17'
$ws.Range("E14").Value = '  +2.08%  '
$ws.Range("D15").Value = '60.737.13'
$ws.Range("E15").Value = '  +2.93%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '2.966.18'
$ws.Range("E16").Value = '  +11.49%  '
$ws.Range("B17").Value = 'Avalanche'
$ws.Range("C17").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D17").Value = '21.38'
$ws.Range("E17").Value = '  +2.22%  '
$ws.Range("E18").Value = '  +0.78%  '
$ws.Range("D19").Value = '349.63'
$ws.Range("E19").Value = '  +0.13%  '
$ws.Range("D20").Value = '4.53'
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("D21").Value = '10.59'
$ws.Range("E21").Value = '  +2.27%  '
$ws.Range("E22").Value = '  +2.60%  '
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D24").Value = '63.72'
$ws.Range("E24").Value = '  +2.97%  '
$ws.Range("E25").Value = '  +0.61%  '
$ws.Range("D26").Value = '0.171'
$ws.Range("E26").Value = '  +5.54%  '
$ws.Range("D27").Value = '0.993'
$ws.Range("E27").Value = '  +0.13%  '
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").Value = '7.39'
$ws.Range("E28").Value = '  +3.79%  '
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").Value = '0.0₃0820'
$ws.Range("E29").Value = '  +2.07%  '
$ws.Range("E30").Value = '  +9.38%  '
$ws.Range("E31").Value = '  +0.07%  '
$ws.Range("D32").Value = '19.27'
$ws.Range("E32").Value = '  +1.38%  '
$ws.Range("E33").Value = '  +1.61%  '
$ws.Range("D34").Value = '149.76'
$ws.Range("E34").Value = '  -0.19%  '
$ws.Range("E35").Value = '  +6.70%  '
$ws.Range("D36").Value = '1.25'
$ws.Range("E36").Value = '  +10.07%  '
$ws.Range("D37").Value = '0.953'
$ws.Range("E37").Value = '  -1.68%  '
$ws.Range("D38").Value = '0.883'
$ws.Range("E38").Value = '  +4.99%  '
$ws.Range("E39").Value = '  +8.57%  '
$ws.Range("D40").Value = '36.99'
$ws.Range("E40").Value = '  +0.93%  '
$ws.Range("D41").Value = '3.67'
$ws.Range("E41").Value = '  -1.05%  '
$ws.Range("D42").Value = '283.93'
$ws.Range("E42").Value = '  +1.79%  '
$ws.Range("D43").Value = '20.18'
$ws.Range("E43").Value = '  +2.80%  '
$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").Value = '0.614'
$ws.Range("E44").Value = '  +0.29%  '
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").Value = '0.0992'
$ws.Range("E45").Value = '  +0.46%  '
$ws.Range("D46").Value = '2.147.36'
$ws.Range("E46").Value = '  +8.14%  '
$ws.Range("D47").Value = '0.996'
$ws.Range("E47").Value = '  +0.09%  '
$ws.Range("E48").Value = '  +5.31%  '
$ws.Range("D49").Value = '0.0541'
$ws.Range("E49").Value = '  +2.07%  '
$ws.Range("E50").Value = '  +2.36%  '
$ws.Range("E51").Value = '  +1.62%  '
